$wb = $excel.ActiveWorkbook

# The worksheet that holds the paid-leave data is the sheet named "Paid_leave"
$ws = $wb.Worksheets.Item("Paid_leave")

# Append the three new rows of data (employee_id, date) below the existing
# table, continuing the pattern already present in the sheet. The "date"
# column stores plain text such as "02.10.2020" (not a real Excel date), so
# a leading apostrophe is used to force each value to be kept as literal
# text instead of being auto-converted into a date serial number; the style
# is reset back to Normal afterwards so no quote-prefix formatting lingers.
$ws.Cells.Item(27, 1).Value = 3
$ws.Cells.Item(27, 2).Value = "'11.12.2020"
$ws.Cells.Item(27, 2).Style = "Normal"

$ws.Cells.Item(28, 1).Value = 1
$ws.Cells.Item(28, 2).Value = "'02.12.2020"
$ws.Cells.Item(28, 2).Style = "Normal"

$ws.Cells.Item(29, 1).Value = 2
$ws.Cells.Item(29, 2).Value = "'18.12.2020"
$ws.Cells.Item(29, 2).Style = "Normal"

# Make this sheet active and bring the new rows into view, matching the
# scrolled/selected state captured in the saved workbook.
$ws.Activate()
$excel.Goto($ws.Range("A15"), $false)
$ws.Range("B29").Select()
